$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 2 (shifts existing rows 2-9 down to 3-10)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the TOC filter header/values
$ws.Range("A2").Value = "TOC Filter"
$ws.Range("B2:G2").Value = "All TOCs"

# Update selection to H3 as in the target
$ws.Range("H3").Select()
